# Generate Report for Archive
#
# The localization-status report was regenerated: the entry for
# "f7f33517-d9e8-45af-9297-9d38d63b0f4c.md" moved from "Ready for
# handoff" back to "In Translation". Because the report sheets are
# sorted by status (In Translation rows first, then Ready for handoff
# rows), that file's row now sorts ahead of
# "0f3411bb-ab2d-41fd-90de-a8104a37b160.md" 's row, so the two rows
# (4 and 5) trade places on every sheet (Overview, zh-cn, de-de): the
# file-identifying columns (name / handoff file / handoff datetime)
# swap between row 4 and row 5, while the Status column is set to the
# correct value for whichever file now occupies that row. Hyperlink
# display text is kept in sync with the cell it now labels (the
# hyperlink target addresses and r:id assignments stay anchored to
# their original cell position, matching the unchanged .rels files).

function Swap-CellValue($ws, $addr1, $addr2) {
    $v1 = $ws.Range($addr1).Value2
    $v2 = $ws.Range($addr2).Value2
    $ws.Range($addr1).Value2 = $v2
    $ws.Range($addr2).Value2 = $v1
}

function Get-HyperlinkForCell($ws, $addr) {
    $target = $ws.Range($addr).Address()
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $target) {
            return $hl
        }
    }
    return $null
}

function Swap-HyperlinkDisplay($ws, $addr1, $addr2) {
    $hl1 = Get-HyperlinkForCell $ws $addr1
    $hl2 = Get-HyperlinkForCell $ws $addr2
    if (($hl1 -ne $null) -and ($hl2 -ne $null)) {
        $t1 = $hl1.TextToDisplay
        $t2 = $hl2.TextToDisplay
        $hl1.TextToDisplay = $t2
        $hl2.TextToDisplay = $t1
    }
}

$wb = $excel.ActiveWorkbook

# --- Overview sheet: columns A (file name), B (zh-cn status), C (de-de status) ---
$ws = $wb.Worksheets.Item("Overview")
Swap-CellValue $ws "A4" "A5"
$ws.Range("B4").Value2 = "In Translation"
$ws.Range("C4").Value2 = "In Translation"
$ws.Range("B5").Value2 = "Ready for handoff"
$ws.Range("C5").Value2 = "Ready for handoff"
Swap-HyperlinkDisplay $ws "A4" "A5"

# --- zh-cn sheet: columns A (file name), B (status), C (handoff xlf), D (handoff datetime) ---
$ws = $wb.Worksheets.Item("zh-cn")
Swap-CellValue $ws "A4" "A5"
Swap-CellValue $ws "C4" "C5"
Swap-CellValue $ws "D4" "D5"
$ws.Range("B4").Value2 = "In Translation"
$ws.Range("B5").Value2 = "Ready for handoff"
Swap-HyperlinkDisplay $ws "A4" "A5"
Swap-HyperlinkDisplay $ws "C4" "C5"

# --- de-de sheet: columns A (file name), B (status), C (handoff xlf), D (handoff datetime) ---
$ws = $wb.Worksheets.Item("de-de")
Swap-CellValue $ws "A4" "A5"
Swap-CellValue $ws "C4" "C5"
Swap-CellValue $ws "D4" "D5"
$ws.Range("B4").Value2 = "In Translation"
$ws.Range("B5").Value2 = "Ready for handoff"
Swap-HyperlinkDisplay $ws "A4" "A5"
Swap-HyperlinkDisplay $ws "C4" "C5"

Write-Output "localization-status report rows 4/5 updated on all sheets"
